# Updates the "cryptos" worksheet with refreshed price / volume(1h) figures
# (and restores the correct Kaspa / ImmutableX row ordering), matching the
# GitHub Actions data-refresh commit.
#
# Numeric-looking Price values (single decimal point, e.g. "607.62") are
# written with a leading apostrophe so Excel keeps them as literal text
# instead of silently converting them to numbers (which would drop
# trailing zeros / change formatting, e.g. "1.00" -> 1). Values that are
# not numeric-looking (thousands-dotted prices like "73.030.67", the
# Volume(1h) percentage strings, coin names and links) are assigned as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = '73.030.67'
$ws.Range("E2").Value = '  +5.97%  '

# Row 3 - Ethereum
$ws.Range("D3").Value = '2.659.93'
$ws.Range("E3").Value = '  +6.65%  '

# Row 4 - TetherUSD
$ws.Range("E4").Value = '  -0.03%  '

# Row 5 - BNB
$ws.Range("D5").Value = "'607.62"
$ws.Range("E5").Value = '  +2.83%  '

# Row 6 - Solana
$ws.Range("D6").Value = "'181.30"
$ws.Range("E6").Value = '  +3.69%  '

# Row 7 - USDC
$ws.Range("E7").Value = '  -0.01%  '

# Row 8 - XRP
$ws.Range("D8").Value = "'0.530"
$ws.Range("E8").Value = '  +3.13%  '

# Row 9 - Dogecoin
$ws.Range("D9").Value = "'0.175"
$ws.Range("E9").Value = '  +15.24%  '

# Row 10 - LidoStakedEther
$ws.Range("D10").Value = '2.658.21'
$ws.Range("E10").Value = '  +6.54%  '

# Row 11 - TRON
$ws.Range("E11").Value = '  +1.26%  '

# Row 12 - Cardano
$ws.Range("D12").Value = "'0.354"
$ws.Range("E12").Value = '  +5.59%  '

# Row 13 - Toncoin
$ws.Range("E13").Value = '  +3.03%  '

# Row 14 - ShibaInu
$ws.Range("D14").Value = "'0.0000193"
$ws.Range("E14").Value = '  +11.36%  '

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = '3.134.11'
$ws.Range("E15").Value = '  +5.88%  '

# Row 16 - Avalanche
$ws.Range("D16").Value = "'27.05"
$ws.Range("E16").Value = '  +5.97%  '

# Row 17 - WrappedBTC
$ws.Range("D17").Value = '72.923.58'
$ws.Range("E17").Value = '  +6.04%  '

# Row 18 - WrappedEther
$ws.Range("D18").Value = '2.658.39'
$ws.Range("E18").Value = '  +6.64%  '

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "'385.56"
$ws.Range("E19").Value = '  +7.43%  '

# Row 20 - Chainlink
$ws.Range("D20").Value = "'11.64"
$ws.Range("E20").Value = '  +7.54%  '

# Row 21 - Uniswap
$ws.Range("D21").Value = "'7.94"
$ws.Range("E21").Value = '  +5.90%  '

# Row 22 - Polkadot
$ws.Range("E22").Value = '  +5.75%  '

# Row 23 - SuiNetwork
$ws.Range("D23").Value = "'2.04"
$ws.Range("E23").Value = '  +23.53%  '

# Row 24 - Litecoin
$ws.Range("D24").Value = "'73.63"
$ws.Range("E24").Value = '  +5.49%  '

# Row 25 - NEARProtocol
$ws.Range("E25").Value = '  +7.34%  '

# Row 26 - Dai
$ws.Range("E26").Value = '  +0.20%  '

# Row 27 - Aptos
$ws.Range("D27").Value = "'9.98"
$ws.Range("E27").Value = '  +12.34%  '

# Row 28 - WrappedeETH
$ws.Range("D28").Value = '2.796.63'
$ws.Range("E28").Value = '  +6.51%  '

# Row 29 - Binance-PegBSC-USD
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = '  +1.55%  '

# Row 30 - PEPE
$ws.Range("D30").Value = '0.0₃0977'
$ws.Range("E30").Value = '  +11.99%  '

# Row 31 - Bittensor
$ws.Range("D31").Value = "'540.38"
$ws.Range("E31").Value = '  +6.94%  '

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").Value = "'8.12"
$ws.Range("E32").Value = '  +5.98%  '

# Row 33 - Fetch.AI
$ws.Range("E33").Value = '  +12.12%  '

# Row 34 - PancakeSwap
$ws.Range("E34").Value = '  +5.02%  '

# Row 35 - FirstDigitalUSD
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = '  -0.20%  '

# Row 36 - Monero
$ws.Range("D36").Value = "'162.52"
$ws.Range("E36").Value = '  -0.31%  '

# Row 37 - EthereumClassic
$ws.Range("D37").Value = "'19.42"
$ws.Range("E37").Value = '  +4.76%  '

# Row 38 - now Kaspa (was ImmutableX)
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").Value = "'0.114"
$ws.Range("E38").Value = '  -3.31%  '

# Row 39 - now ImmutableX (was Kaspa)
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").Value = "'1.42"
$ws.Range("E39").Value = '  +9.95%  '

# Row 40 - WhiteBITCoin
$ws.Range("E40").Value = '  +2.54%  '

# Row 41 - Stacks
$ws.Range("E41").Value = '  +10.80%  '

# Row 42 - dogwifhat
$ws.Range("E42").Value = '  +17.60%  '

# Row 43 - RenderToken
$ws.Range("D43").Value = "'5.15"
$ws.Range("E43").Value = '  +9.05%  '

# Row 44 - USDe
$ws.Range("E44").Value = '  +0.19%  '

# Row 45 - PolygonEcosystemToken
$ws.Range("D45").Value = "'0.337"
$ws.Range("E45").Value = '  +6.60%  '

# Row 46 - OKB
$ws.Range("D46").Value = "'39.81"
$ws.Range("E46").Value = '  +2.91%  '

# Row 47 - Aave
$ws.Range("D47").Value = "'152.20"
$ws.Range("E47").Value = '  +2.32%  '

# Row 48 - Filecoin
$ws.Range("D48").Value = "'3.71"
$ws.Range("E48").Value = '  +5.28%  '

# Row 49 - ARBITRUM
$ws.Range("E49").Value = '  +7.82%  '

# Row 50 - Optimism
$ws.Range("E50").Value = '  +11.30%  '

# Row 51 - BabyDogeCoin
$ws.Range("E51").Value = '  +11.11%  '
